$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D24").Value = "Current DC MPPT x"
$ws.Range("D25").Value = "Current DC MPPT x Input x"
$ws.Range("D30").Value = "Power DC MPPT x"
$ws.Range("D41").Value = "Voltage DC MPPT x"
$ws.Range("D42").Value = "Voltage DC MPPT x Input x"
